# Add two new date columns (27-sep, 28-sep) to the Dataframe Fam sheet.
# Columns BZ (78) and CA (79) are appended after the existing BY (77) column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new date labels
$ws.Cells.Item(1, 78).Value = "27-sep"
$ws.Cells.Item(1, 79).Value = "28-sep"

function Set-DataCell {
    param($row, $col, $val)
    $c = $ws.Cells.Item($row, $col)
    # Match the existing numeric data-cell formatting (centered integer,
    # same as style used by the rest of the row) before writing the value,
    # so the new cell reuses the workbook's existing style entry.
    $c.HorizontalAlignment = -4108
    $c.NumberFormat = "0"
    $c.Value = $val
}

# Data rows 2-11: new counts for the two new date columns
Set-DataCell 2 78 10
Set-DataCell 2 79 12

Set-DataCell 3 78 16
Set-DataCell 3 79 18

Set-DataCell 4 78 11
Set-DataCell 4 79 11

Set-DataCell 5 78 13
Set-DataCell 5 79 14

Set-DataCell 6 78 8
Set-DataCell 6 79 8

Set-DataCell 7 78 15
Set-DataCell 7 79 16

Set-DataCell 8 78 14
Set-DataCell 8 79 13

Set-DataCell 9 78 16
Set-DataCell 9 79 15

Set-DataCell 10 78 23
Set-DataCell 10 79 31

Set-DataCell 11 78 8
Set-DataCell 11 79 7

# Match the author's final selection recorded in the saved file
$ws.Range("BZ12").Select()
